# Add 2022-Q4 data:
#  1. Insert a new worksheet named "2022-Q4" right after "总计" (i.e. before
#     the sheet that is currently in position 2, "2022-Q3"). All the other
#     quarter sheets keep their name/content and simply shift one tab to the
#     right.
#  2. Populate the new sheet with the fund-holdings table for 2022-Q4.
#  3. Insert a matching row at the top of the "总计" (summary) sheet's data
#     and shift the existing rows down by one, keeping all their values.
#  4. Restore the originally-selected tab ("2020-Q4") as the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before the current 2nd tab.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Fill in the 2022-Q4 holdings table.
# ---------------------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160
$newSheet.Range("B1:H1").Borders.LineStyle = 1

$q4Data = @(
    @(0, "257030", "国联安优势混合",                     "8.55", "89.30", "4.46", "0.3813", 8),
    @(1, "255010", "国联安稳健混合",                     "2.22", "68.97", "3.04", "0.0675", 10),
    @(2, "005576", "华泰柏瑞新金融地产灵活配置混合A",   "0.77", "94.17", "6.61", "0.0509", 5),
    @(3, "001244", "华泰柏瑞量化智慧灵活配置混合A",     "2.89", "93.57", "0.66", "0.0191", 4),
    @(4, "016374", "华泰柏瑞新金融地产灵活配置混合C",   "0.22", "94.17", "6.61", "0.0145", 5),
    @(5, "003760", "国泰中证500指数增强A",               "0.48", "91.42", "1.62", "0.0078", 10),
    @(6, "006104", "华泰柏瑞量化智慧灵活配置混合C",     "0.38", "93.57", "0.66", "0.0025", 4),
    @(7, "003761", "国泰中证500指数增强C",               "0.04", "91.42", "1.62", "0.0006", 10)
)

$row = 2
foreach ($item in $q4Data) {
    $newSheet.Cells.Item($row, 1).Value = $item[0]
    $newSheet.Cells.Item($row, 1).Font.Bold = $true
    $newSheet.Cells.Item($row, 1).HorizontalAlignment = -4108
    $newSheet.Cells.Item($row, 1).VerticalAlignment = -4160
    $newSheet.Cells.Item($row, 1).Borders.LineStyle = 1

    $newSheet.Cells.Item($row, 2).Value = "'" + $item[1]
    $newSheet.Cells.Item($row, 3).Value = $item[2]
    $newSheet.Cells.Item($row, 4).Value = "'" + $item[3]
    $newSheet.Cells.Item($row, 5).Value = "'" + $item[4]
    $newSheet.Cells.Item($row, 6).Value = "'" + $item[5]
    $newSheet.Cells.Item($row, 7).Value = "'" + $item[6]
    $newSheet.Cells.Item($row, 8).Value = $item[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift existing rows down by one and
#    insert the new 2022-Q4 totals at row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q3", 4, 0.08),
    @("2022-Q1", 2, 0.01),
    @("2021-Q3", 2, 0.02),
    @("2021-Q2", 2, 0.28),
    @("2021-Q1", 11, 0.72),
    @("2020-Q4", 6, 0.12)
)

for ($i = $summaryRows.Length - 1; $i -ge 0; $i--) {
    $targetRow = $i + 3
    $vals = $summaryRows[$i]
    $totalSheet.Cells.Item($targetRow, 1).Value = $i + 1
    $totalSheet.Cells.Item($targetRow, 1).Font.Bold = $true
    $totalSheet.Cells.Item($targetRow, 1).HorizontalAlignment = -4108
    $totalSheet.Cells.Item($targetRow, 1).VerticalAlignment = -4160
    $totalSheet.Cells.Item($targetRow, 1).Borders.LineStyle = 1
    $totalSheet.Cells.Item($targetRow, 2).Value = $vals[0]
    $totalSheet.Cells.Item($targetRow, 3).Value = $vals[1]
    $totalSheet.Cells.Item($targetRow, 4).Value = $vals[2]
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 1).Font.Bold = $true
$totalSheet.Cells.Item(2, 1).HorizontalAlignment = -4108
$totalSheet.Cells.Item(2, 1).VerticalAlignment = -4160
$totalSheet.Cells.Item(2, 1).Borders.LineStyle = 1
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.54

# ---------------------------------------------------------------------
# 4. Restore the original active tab ("2020-Q4").
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
